# Add the new "ТДР..." reference numbers into column E (and duplicate a
# couple of them into column H), matching the source workbook update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value  = "ТДР00000275 от 16.04.2021"
$ws.Range("E5").Value  = "ТДР00000329 от 11.05.2021"
$ws.Range("E6").Value  = "ТДР00000384 от 15.06.2021"
$ws.Range("E7").Value  = "ТДР00000464 от 13.07.2021"
$ws.Range("E8").Value  = "ТДР00000641 от 10.08.2021"
$ws.Range("H8").Value  = "ТДР00000766 от 07.10.2021"
$ws.Range("E9").Value  = "ТДР00000667 от 31.08.2021"
$ws.Range("E10").Value = " ТДР00000766 от 07.10.2021"
$ws.Range("E11").Value = "ТДР00000872 от 23.11.2021"
$ws.Range("H11").Value = "ТДР00000872 от 23.11.2021"
$ws.Range("E12").Value = "ТДР00000959 от 15.12.2021"
$ws.Range("E13").Value = "ТДР00000015 от 14.01.2022"

# View state: zoomed to 130%, scrolled so column B is leftmost, H8 selected.
$excel.ActiveWindow.Zoom = 130
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("H8").Select()
